# Apply the StructureDefinition-reimburse-type.xlsx update:
#  - Update IG metadata on the "Metadata" sheet (URL, Version, Date, Publisher)
#  - Clear the "Constraint(s)" value for the top-level Extension row on the "Elements" sheet

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/reimburse-type"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

$elements = $wb.Worksheets.Item("Elements")
$elements.Range("AI2").Value = ""
